$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 with the new "time_taken" column label
$ws.Range("F1").Value = "time_taken"

# Copy formatting from the existing header cell (E1) onto F1
# so the new header matches the bold/bordered header style.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data cells F2:F40 with time_taken timestamps (stored as text)
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:50:51.532628"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:50:51.532639"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:50:51.532643"
$ws.Cells.Item(5, 6).Value = "2021-10-05 10:50:51.532645"
$ws.Cells.Item(6, 6).Value = "2021-10-05 10:50:51.532648"
$ws.Cells.Item(7, 6).Value = "2021-10-05 10:50:51.532651"
$ws.Cells.Item(8, 6).Value = "2021-10-05 10:50:51.532654"
$ws.Cells.Item(9, 6).Value = "2021-10-05 10:50:51.532656"
$ws.Cells.Item(10, 6).Value = "2021-10-05 10:50:51.532659"
$ws.Cells.Item(11, 6).Value = "2021-10-05 10:50:51.532661"
$ws.Cells.Item(12, 6).Value = "2021-10-05 10:50:51.532664"
$ws.Cells.Item(13, 6).Value = "2021-10-05 10:50:51.532667"
$ws.Cells.Item(14, 6).Value = "2021-10-05 10:50:51.532669"
$ws.Cells.Item(15, 6).Value = "2021-10-05 10:50:51.532672"
$ws.Cells.Item(16, 6).Value = "2021-10-05 10:50:51.532674"
$ws.Cells.Item(17, 6).Value = "2021-10-05 10:50:51.532677"
$ws.Cells.Item(18, 6).Value = "2021-10-05 10:50:51.532680"
$ws.Cells.Item(19, 6).Value = "2021-10-05 10:50:51.532682"
$ws.Cells.Item(20, 6).Value = "2021-10-05 10:50:51.532685"
$ws.Cells.Item(21, 6).Value = "2021-10-05 10:50:51.532687"
$ws.Cells.Item(22, 6).Value = "2021-10-05 10:50:51.532690"
$ws.Cells.Item(23, 6).Value = "2021-10-05 10:50:51.532692"
$ws.Cells.Item(24, 6).Value = "2021-10-05 10:50:51.532695"
$ws.Cells.Item(25, 6).Value = "2021-10-05 10:50:51.532697"
$ws.Cells.Item(26, 6).Value = "2021-10-05 10:50:51.532700"
$ws.Cells.Item(27, 6).Value = "2021-10-05 10:50:51.532703"
$ws.Cells.Item(28, 6).Value = "2021-10-05 10:50:51.532705"
$ws.Cells.Item(29, 6).Value = "2021-10-05 10:50:51.532708"
$ws.Cells.Item(30, 6).Value = "2021-10-05 10:50:51.532711"
$ws.Cells.Item(31, 6).Value = "2021-10-05 10:50:51.532713"
$ws.Cells.Item(32, 6).Value = "2021-10-05 10:50:51.532716"
$ws.Cells.Item(33, 6).Value = "2021-10-05 10:50:51.532718"
$ws.Cells.Item(34, 6).Value = "2021-10-05 10:50:51.532721"
$ws.Cells.Item(35, 6).Value = "2021-10-05 10:50:51.532724"
$ws.Cells.Item(36, 6).Value = "2021-10-05 10:50:51.532726"
$ws.Cells.Item(37, 6).Value = "2021-10-05 10:50:51.532728"
$ws.Cells.Item(38, 6).Value = "2021-10-05 10:50:51.532731"
$ws.Cells.Item(39, 6).Value = "2021-10-05 10:50:51.532734"
$ws.Cells.Item(40, 6).Value = "2021-10-05 10:50:51.532736"
